$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style to use for date/time columns (B, C, H) - copy from an existing date cell
$dateFormat = $ws.Cells.Item(2, 2).NumberFormat

# ---- Row 5 ----
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 44523.82418981481
$ws.Cells.Item(5, 2).NumberFormat = $dateFormat
$ws.Cells.Item(5, 3).Value = 44523.8263425926
$ws.Cells.Item(5, 3).NumberFormat = $dateFormat
$ws.Cells.Item(5, 4).Value = "IP Address"
$ws.Cells.Item(5, 5).Value = 100
$ws.Cells.Item(5, 6).Value = 186
$ws.Cells.Item(5, 7).Value = $true
$ws.Cells.Item(5, 8).Value = 44523.82635416667
$ws.Cells.Item(5, 8).NumberFormat = $dateFormat
$ws.Cells.Item(5, 9).Value = "1dabec"
$ws.Cells.Item(5, 10).Value = "ebola %>%`n  pivot_longer(``Cases_Guinea``:last_col()) %>%`n  separate(name, into = c(""case_death"", ""country""), sep = ""_"") %>%`n  drop_na()"

# ---- Row 6 ----
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 44523.82447916667
$ws.Cells.Item(6, 2).NumberFormat = $dateFormat
$ws.Cells.Item(6, 3).Value = 44523.8337962963
$ws.Cells.Item(6, 3).NumberFormat = $dateFormat
$ws.Cells.Item(6, 4).Value = "IP Address"
$ws.Cells.Item(6, 5).Value = 100
$ws.Cells.Item(6, 6).Value = 804
$ws.Cells.Item(6, 7).Value = $true
$ws.Cells.Item(6, 8).Value = 44523.8337962963
$ws.Cells.Item(6, 8).NumberFormat = $dateFormat
$ws.Cells.Item(6, 9).Value = "1davec"
$ws.Cells.Item(6, 10).Value = "ebola <- read_csv(""ebola.csv"")`nebola_tidy <- ebola %>%`n  pivot_longer(``Cases_Guinea``:last_col()) %>%`n  separate(name, into = c(""case_death"", ""country""), sep = ""_"") %>% drop_na()"

# ---- Row 7 ----
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 44523.18385416667
$ws.Cells.Item(7, 2).NumberFormat = $dateFormat
$ws.Cells.Item(7, 3).Value = 44523.18526620371
$ws.Cells.Item(7, 3).NumberFormat = $dateFormat
$ws.Cells.Item(7, 4).Value = "Spam"
$ws.Cells.Item(7, 5).Value = 50
$ws.Cells.Item(7, 6).Value = 121
$ws.Cells.Item(7, 7).Value = $false
$ws.Cells.Item(7, 8).Value = 44523.87954861111
$ws.Cells.Item(7, 8).NumberFormat = $dateFormat
$ws.Cells.Item(7, 9).Value = "2nesch"

Write-Host "rows added"
